$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E (shared string "PresskPa")
$ws.Range("E2").Value = "PresskPa"

# New column E: PressBar (column D) converted to kPa (x100)
# Row 3 gets its own (non-shared) formula, matching how column C's
# formula chain (C3 alone, then C4:C26 shared) was originally built.
$ws.Range("E3").Formula = "=D3*100"
$ws.Range("E4:E26").Formula = "=D4*100"

# Leave the selection where the edit finished, as captured in the workbook
$ws.Range("J22").Select()
